# Slide 2, "Content Placeholder 11" shape: two small wording tweaks per the
# commit ("Made very minor change on slide 2 as bugs seem to be resolved.")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# 1) Re-word the "code is ~100% done" status line (paragraph 2).
$statusPara = $tr.Paragraphs(2)
$statusPara.Text = "Code appears to be at 100pct as all minor bugs seem to be fixed"

# 2) Andre Nasrah's contribution paragraph (paragraph 9) was split across two
#    runs - "...coordination " and "of meetings." - apparently from a
#    trailing edit. Merge them back into a single run reading "...
#    coordination of meetings."
$andrePara = $tr.Paragraphs(9)
$secondRun = $andrePara.Runs(2)
$secondRun.Text = ""
$secondRun.Delete()

$andrePara = $tr.Paragraphs(9)
$firstRun = $andrePara.Runs(1)
$firstRun.Text = "Andre Nasrah contributed conceptual " + [char]0x2013 + " stripped down code before as well implementing converterclass.java to have a class object associated with program, basic layout in activity_main.xml, and basis of documentation.  Contributed to Github regarding bug editing. Provided test cases along with updates which were regularly pushed to Github. Zoom meeting setups and coordination of meetings."

Write-Host "slide 2 text updated"
